$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subjectGmail = " 🔎 : Prueba automatizacion Gmail/whatasapp "
$fromJoan     = " 📭 : Joan Martinez <joan_martinez.olivares@hotmail.com>"
$bodyGmail    = "Prueba automatizacion Gmail/whatasapp "
$newCol       = "📩 NUEVO 📩"

$subjectRenata = " 🔎 : Prueba con renata"
$bodyRenata    = "Prueba con renata "

# Rows 25-39: 15 rows with the "Prueba automatizacion Gmail/whatasapp" content
for ($r = 25; $r -le 39; $r++) {
    $ws.Cells.Item($r, 1).Value = $subjectGmail
    $ws.Cells.Item($r, 2).Value = $fromJoan
    $ws.Cells.Item($r, 3).Value = $bodyGmail
    $ws.Cells.Item($r, 4).Value = $newCol
}

# Rows 40-41: 2 rows with the "Prueba con renata" content
for ($r = 40; $r -le 41; $r++) {
    $ws.Cells.Item($r, 1).Value = $subjectRenata
    $ws.Cells.Item($r, 2).Value = $fromJoan
    $ws.Cells.Item($r, 3).Value = $bodyRenata
    $ws.Cells.Item($r, 4).Value = $newCol
}
